$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 20:39"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5270736
$ws.Range("C4").Value = 19290
$ws.Range("D4").Value = 2719695
$ws.Range("E4").Value = 2384170
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 679
$ws.Range("H4").Value = 166871

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2328405
$ws.Range("C6").Value = 61252
$ws.Range("D6").Value = 1636614
$ws.Range("E6").Value = 645603
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 835
$ws.Range("H6").Value = 46188

# Row 15: Reino Unido
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 312789
$ws.Range("C15").Value = 1148
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 46526

# Row 23: Francia
$ws.Range("A23").Value = "Francia"
$ws.Range("B23").Value = 204172
$ws.Range("C23").Value = 1397
$ws.Range("D23").Value = 82836
$ws.Range("E23").Value = 90996
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 30340

# Row 31: Ecuador
$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 95563
$ws.Range("C31").Value = 862
$ws.Range("D31").Value = 78610
$ws.Range("E31").Value = 11002
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 19
$ws.Range("H31").Value = 5951

# Row 33: Israel
$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 86147
$ws.Range("C33").Value = 1425
$ws.Range("D33").Value = 60055
$ws.Range("E33").Value = 25470
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 622

# Row 60: Marruecos
$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 35195
$ws.Range("C60").Value = 1132
$ws.Range("D60").Value = 25385
$ws.Range("E60").Value = 9277
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 17
$ws.Range("H60").Value = 533

# Row 98: Libano
$ws.Range("A98").Value = "Libano"
$ws.Range("B98").Value = 7121
$ws.Range("C98").Value = 309
$ws.Range("D98").Value = 2377
$ws.Range("E98").Value = 4657
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 7
$ws.Range("H98").Value = 87

# Row 107: Zimbabue
$ws.Range("A107").Value = "Zimbabue"
$ws.Range("B107").Value = 4818
$ws.Range("C107").Value = 70
$ws.Range("D107").Value = 1524
$ws.Range("E107").Value = 3190
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 104

# Row 110: Republica de Africa Central
$ws.Range("A110").Value = "Republica de Africa Central"
$ws.Range("B110").Value = 4645
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 1723
$ws.Range("E110").Value = 2861
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 61

# Row 119: Cuba
$ws.Range("A119").Value = "Cuba"
$ws.Range("B119").Value = 3093
$ws.Range("C119").Value = 47
$ws.Range("D119").Value = 2472
$ws.Range("E119").Value = 533
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 88

# Row 120: Mayotte
$ws.Range("A120").Value = "Mayotte"
$ws.Range("B120").Value = 3068
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 2835
$ws.Range("E120").Value = 194
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 39

# Row 142: Siria
$ws.Range("A142").Value = "Siria"
$ws.Range("B142").Value = 1327
$ws.Range("C142").Value = 72
$ws.Range("D142").Value = 385
$ws.Range("E142").Value = 889
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 53

# Row 143: Uganda
$ws.Range("A143").Value = "Uganda"
$ws.Range("B143").Value = 1313
$ws.Range("C143").Value = 16
$ws.Range("D143").Value = 1138
$ws.Range("E143").Value = 166
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 9

# Row 144: Letonia
$ws.Range("A144").Value = "Letonia"
$ws.Range("B144").Value = 1293
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 1078
$ws.Range("E144").Value = 183
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 32

# Row 145: Jordania
$ws.Range("A145").Value = "Jordania"
$ws.Range("B145").Value = 1283
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 1189
$ws.Range("E145").Value = 83
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 11

# Row 146: Georgia
$ws.Range("A146").Value = "Georgia"
$ws.Range("B146").Value = 1264
$ws.Range("C146").Value = 14
$ws.Range("D146").Value = 1054
$ws.Range("E146").Value = 193
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 17

# Row 148: Liberia
$ws.Range("A148").Value = "Liberia"
$ws.Range("B148").Value = 1250
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 736
$ws.Range("E148").Value = 433
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 81

# Row 149: Burkina Faso
$ws.Range("A149").Value = "Burkina Faso"
$ws.Range("B149").Value = 1211
$ws.Range("C149").Value = 7
$ws.Range("D149").Value = 990
$ws.Range("E149").Value = 167
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 54

# Row 213: Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
